$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("breakdown") to hold the "country" field.
$ws.Columns("B").Insert()

# Header + per-row country values (new column B), keyed by row number.
$countries = @{
    1  = "country"
    2  = "Spain"
    3  = "Spain"
    4  = "Netherlands"
    5  = "Denmark"
    6  = "Denmark"
    7  = "Switzerland"
    8  = "Switzerland"
    9  = "Brazil"
    10 = "Brazil"
    11 = "New York"
    12 = "Italy"
    13 = "Italy"
    14 = "China, Wuhan"
    15 = "Luxembourg"
    16 = "England"
    17 = "England"
    18 = "Switzerland"
    19 = "Denmark"
    20 = "Spain"
    21 = "England"
    22 = "Switzerland"
    23 = "Switzerland"
    24 = "New York"
}

foreach ($row in $countries.Keys) {
    $ws.Cells.Item($row, 2).Value = $countries[$row]
}

# Match the author's post-edit selection.
$ws.Range("B25").Select()
